$p = $ppt.ActivePresentation

# --- Update the date placeholder text (8/29/2017 -> 9/11/2017) across the
#     slide master and every custom layout that carries a Date Placeholder. ---
function Update-DatePlaceholders($shapes) {
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $shape = $shapes.Item($j)
        if ($shape.Name -like "Date Placeholder*" -and $shape.HasTextFrame) {
            $tr = $shape.TextFrame.TextRange
            if ($tr.Text -eq "8/29/2017") {
                $tr.Text = "9/11/2017"
            }
        }
    }
}

Update-DatePlaceholders $p.SlideMaster.Shapes

for ($i = 1; $i -le $p.SlideMaster.CustomLayouts.Count; $i++) {
    $layout = $p.SlideMaster.CustomLayouts.Item($i)
    Update-DatePlaceholders $layout.Shapes
}

# --- Bump the version string on the binder cover slide. ---
$slide1 = $p.Slides.Item(1)
for ($j = 1; $j -le $slide1.Shapes.Count; $j++) {
    $shape = $slide1.Shapes.Item($j)
    if ($shape.HasTextFrame -and $shape.TextFrame.TextRange.Text -eq "Version 1.4") {
        $shape.TextFrame.TextRange.Text = "Version 2.0"
    }
}
